# Scen_TRA_Max_Growth.xlsx
# "Modify seed value for UC growth rate and modify AF for private cars"
#
# The UC-growth-rate tables in rows 9/10/12/13/15/16 referenced the 2025
# (column J) / "10 years out" (column L) growth figures for LGV/MGV/HGV.
# The seed year for those max-growth-rate lookups moves back to 2020
# (column E) instead. Rows 8/11/14/17 (the "AF"/FCV variants for
# Cars/LGV/MGV/HGV) are changed the same way, and the red "seed" highlight
# column in the underlying growth table (rows 42-46) moves from column J
# to column E to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update the seed-year growth-rate formulas (L9:L17) to reference the
#    2020 column (E) of the HGV/MGV/LGV/PC growth table (rows 43-46)
#    instead of the old 2025/L-column references.
# ---------------------------------------------------------------------
$ws.Range("L9").Formula  = "=-E45/1000"
$ws.Range("L10").Formula = "=-E45/1000"
$ws.Range("L11").Formula = "=-E45/1000"

$ws.Range("L12").Formula = "=-E44/1000"
$ws.Range("L13").Formula = "=-E44/1000"
$ws.Range("L14").Formula = "=-E44/1000"

$ws.Range("L15").Formula = "=-E43/1000"
$ws.Range("L16").Formula = "=-E43/1000"
$ws.Range("L17").Formula = "=-E43/1000"

# ---------------------------------------------------------------------
# 2) Move the red "seed" highlight format from column J (2025) to
#    column E (2020) in the growth table, rows 42-46.
# ---------------------------------------------------------------------
$ws.Range("J43").Copy() | Out-Null
$ws.Range("E43:E46").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("D43").Copy() | Out-Null
$ws.Range("J42:J46").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Move the active selection from L13 to L28 and scroll the sheet back
#    to the top (clear the old "topLeftCell = A10" frozen scroll position).
# ---------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$ws.Range("L28").Select() | Out-Null
